$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1130968.1
$ws.Range("J17").Value = 1130968.1
$ws.Range("L17").Value = 3392904.3
$ws.Range("N17").Value = -3393240.3

$ws.Range("H76").Value = 76927100
$ws.Range("I76").Value = 100002880
$ws.Range("J76").Value = 7836
$ws.Range("K76").Value = 100002880
$ws.Range("L76").Value = 7836
$ws.Range("M76").Value = -100002565
$ws.Range("N76").Value = -8466

$ws.Range("H79").Value = 76927100
$ws.Range("I79").Value = 100002880
$ws.Range("J79").Value = 7836
$ws.Range("K79").Value = 100002880
$ws.Range("L79").Value = 7836
$ws.Range("M79").Value = -100001788
$ws.Range("N79").Value = -10020

$ws.Range("H116").Value = 5929.811
$ws.Range("I116").Value = 9970.076999999999
$ws.Range("J116").Value = 3741.3333
$ws.Range("K116").Value = 9970.076999999999
$ws.Range("L116").Value = 3741.3333
$ws.Range("M116").Value = -6528.076999999999
$ws.Range("N116").Value = -10625.3333

$ws.Range("H129").Value = 1634.1086
$ws.Range("I129").Value = 398.14285
$ws.Range("J129").Value = 2174.8438
$ws.Range("K129").Value = 1194.42855
$ws.Range("L129").Value = 6524.5314
$ws.Range("M129").Value = 3805.57145
$ws.Range("N129").Value = -16524.5314

$ws.Range("H132").Value = 211391.62
$ws.Range("I132").Value = 3135.718
$ws.Range("J132").Value = 1113833.9
$ws.Range("K132").Value = 9407.153999999999
$ws.Range("L132").Value = 3341501.7
$ws.Range("M132").Value = -6877.153999999999
$ws.Range("N132").Value = -3346561.7

$ws.Range("H139").Value = 58500
$ws.Range("J139").Value = 58500
$ws.Range("L139").Value = 58500
$ws.Range("N139").Value = -68780

$ws.Range("H140").Value = 65780
$ws.Range("J140").Value = 65780
$ws.Range("L140").Value = 65780
$ws.Range("N140").Value = -76140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1081.25
$ws.Range("I61").Value = 930
$ws.Range("J61").Value = 1333.3334
$ws.Range("K61").Value = 930
$ws.Range("L61").Value = 1333.3334
$ws.Range("M61").Value = -718
$ws.Range("N61").Value = -1757.3334

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 69500
$ws.Range("J135").Value = 69500
$ws.Range("L135").Value = 69500
$ws.Range("N135").Value = -79640

$ws.Range("H136").Value = 1081.25
$ws.Range("I136").Value = 930
$ws.Range("J136").Value = 1333.3334
$ws.Range("K136").Value = 2790
$ws.Range("L136").Value = 4000.0002
$ws.Range("M136").Value = -240
$ws.Range("N136").Value = -9100.0002

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26637.56
$ws.Range("I134").Value = 1192.8077
$ws.Range("J134").Value = 70741.8
$ws.Range("K134").Value = 3578.4231
$ws.Range("L134").Value = 212225.4
$ws.Range("M134").Value = -1043.4231
$ws.Range("N134").Value = -217295.4

$ws.Range("H140").Value = 55780
$ws.Range("J140").Value = 55780
$ws.Range("L140").Value = 55780
$ws.Range("N140").Value = -66140

$ws.Range("H141").Value = 36500
$ws.Range("J141").Value = 36500
$ws.Range("L141").Value = 36500
$ws.Range("N141").Value = -46860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3446.6875
$ws.Range("I99").Value = 2962.25
$ws.Range("J99").Value = 4900
$ws.Range("K99").Value = 2962.25
$ws.Range("L99").Value = 4900
$ws.Range("M99").Value = -1464.25
$ws.Range("N99").Value = -7896

$ws.Range("H107").Value = 312.17242
$ws.Range("I107").Value = 226.22728
$ws.Range("J107").Value = 582.2857
$ws.Range("K107").Value = 226.22728
$ws.Range("L107").Value = 582.2857
$ws.Range("M107").Value = 1693.77272
$ws.Range("N107").Value = -4422.2857

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 3446.6875
$ws.Range("I126").Value = 2962.25
$ws.Range("J126").Value = 4900
$ws.Range("K126").Value = 8886.75
$ws.Range("L126").Value = 14700
$ws.Range("M126").Value = -6416.75
$ws.Range("N126").Value = -19640

$ws.Range("H132").Value = 3759.9
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 4619.8
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 13859.4
$ws.Range("M132").Value = -6170
$ws.Range("N132").Value = -18919.4

$ws.Range("H134").Value = 17859536
$ws.Range("I134").Value = 2110.3333
$ws.Range("K134").Value = 6330.999899999999
$ws.Range("M134").Value = -3795.999899999999

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 929.55554
$ws.Range("I132").Value = 776.8889
$ws.Range("J132").Value = 1082.2222
$ws.Range("K132").Value = 6992.0001
$ws.Range("L132").Value = 9739.9998
$ws.Range("M132").Value = -4462.0001
$ws.Range("N132").Value = -14799.9998

$ws.Range("H133").Value = 1430
$ws.Range("I133").Value = 1430
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 4290
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 770
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 4568.1953
$ws.Range("I134").Value = 1482.1538
$ws.Range("J134").Value = 6001
$ws.Range("K134").Value = 4446.4614
$ws.Range("L134").Value = 18003
$ws.Range("M134").Value = 623.5385999999999
$ws.Range("N134").Value = -28143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20462.887
$ws.Range("I132").Value = 1189.2051
$ws.Range("J132").Value = 74153.86
$ws.Range("K132").Value = 3567.615299999999
$ws.Range("L132").Value = 222461.58
$ws.Range("M132").Value = -1037.615299999999
$ws.Range("N132").Value = -227521.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 39938
$ws.Range("J43").Value = 39938
$ws.Range("L43").Value = 39938
$ws.Range("N43").Value = -40324

$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990

$ws.Range("H122").Value = 2481.25
$ws.Range("I122").Value = 2030
$ws.Range("J122").Value = 3233.3333
$ws.Range("K122").Value = 6090
$ws.Range("L122").Value = 9699.999899999999
$ws.Range("M122").Value = -3640
$ws.Range("N122").Value = -14599.9999

$ws.Range("H141").Value = 53873.332
$ws.Range("J141").Value = 53873.332
$ws.Range("L141").Value = 53873.332
$ws.Range("N141").Value = -64233.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 43266
$ws.Range("J125").Value = 43266
$ws.Range("L125").Value = 43266
$ws.Range("N125").Value = -53106

$ws.Range("H132").Value = 2846.3333
$ws.Range("I132").Value = 908.825
$ws.Range("J132").Value = 7405.1763
$ws.Range("K132").Value = 2726.475
$ws.Range("L132").Value = 22215.5289
$ws.Range("M132").Value = -196.4750000000004
$ws.Range("N132").Value = -27275.5289
